$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing column C. This shifts the old
# "Другой язык ..." column (C) to D, and the old "Расширенный комментарий"
# column (D) to E.
$ws.Columns("C").Insert()

# The header wording lost its " | обработка" suffix on both language columns.
$ws.Range("B1").Value = "Переводимый язык (основной диалект) [система письма]"

# New column C repeats the "translated language" header and holds the
# wordform breakdown (django_toolbar wordform parsing columns).
$ws.Range("C1").Value = "Переводимый язык (основной диалект) [система письма]"
$ws.Range("C2").Value = "словоформаА1 | словоформаА2"
$ws.Range("C3").Value = "словоформаБ1"

$ws.Range("D1").Value = "Другой язык (основной диалект) [система письма]"

# Header row is shorter now that it only holds single-line headers.
$ws.Rows("1").RowHeight = 30

# Selection moves along with the content that used to sit in C2.
$ws.Range("D2").Select()
